$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(6,20,20,20,1,2.5,1,0.5,0,15,3,0,0,0,0.01,0.1,-0.4726188122638625,17.81598639631854,-11.06931823895891,0,"blink"),
    @(6,20,20,20,1,2.5,1,0.5,0,15,3,0,0,0,0.01,0.1,-0.6302600461921624,15.34914611787046,-18.50947527611086,1,"templerun"),
    @(6,20,20,20,1,2.5,1,0.5,0,15,3,0,0,0,0.01,0.1,-0.7159122124080665,8.47388609802706,-15.28942257020605,2,"sudoku"),
    @(6,20,20,20,1,2.5,1,0.5,0,15,3,0,0,0,0.01,0.1,-1.07867557886635,6.789725143791687,-20.90523570871855,3,"flow"),
    @(6,20,20,20,1,2.5,1,0.5,0,15,3,0,0,0,0.01,0.01,-1.22771535155418,5.372075960778282,-1.501134687158749,4,"wordsearch")
)

$rowIndex = 2
foreach ($row in $data) {
    $colIndex = 1
    foreach ($val in $row) {
        $ws.Cells.Item($rowIndex, $colIndex).Value = $val
        $colIndex++
    }
    $rowIndex++
}
